$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7095779180526733
$ws.Range("B1").Value = 2.765387296676636
$ws.Range("C1").Value = 3.949029207229614
$ws.Range("D1").Value = 1.372920393943787
$ws.Range("E1").Value = 0.7773333191871643
